$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix capitalization of names in column C (Họ và Tên) ---
$ws.Range("C2").Value  = "Đặng Văn Cương"
$ws.Range("C4").Value  = "Vũ Văn Sướng"
$ws.Range("C6").Value  = "Bạch Văn Toàn"
$ws.Range("C9").Value  = "Bạch Văn Trọng"
$ws.Range("C10").Value = "Bạch Văn Tuấn"
$ws.Range("C13").Value = "Đặng Thị Hương"
$ws.Range("C14").Value = "Bạch Thuỳ Linh"
$ws.Range("C16").Value = "Vũ Thuỳ Linh"
$ws.Range("C23").Value = "Đặng Đình Khải"
$ws.Range("C26").Value = "Đặng Thị Mây"
$ws.Range("C28").Value = "Bạch Văn Thuỷ"
$ws.Range("C35").Value = "Đặng Khánh Hưng"

# --- View changes: zoom to 130% and move the selection ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("J33").Select()

# --- Conditional formatting: highlight duplicate values in J2:J16 ---
# (Excel's built-in "Highlight Duplicate Values" preset: dark red text on
# a light red fill.) Two identical rule objects are created and one is
# removed, matching the leftover unused dxf entry Excel leaves behind in
# styles.xml when a duplicate-values rule is edited/reapplied via the UI.
$rng = $ws.Range("J2:J16")

$fc1 = $rng.FormatConditions.AddUniqueValues()
$fc1.DupeUnique = 1
$fc1.Font.Color = 393372
$fc1.Interior.Color = 13551615

$fc2 = $rng.FormatConditions.AddUniqueValues()
$fc2.DupeUnique = 1
$fc2.Font.Color = 393372
$fc2.Interior.Color = 13551615
$fc2.Delete()
